$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.78343898000435
$ws.Range("C2").Value = 8.320682808848133
$ws.Range("D2").Value = 3.944530293171527
$ws.Range("E2").Value = 11.80144205942174
$ws.Range("F2").Value = 21.68843984460188
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("M2").Value = 14.2283810455674
$ws.Range("N2").Value = 16.55318618368921
$ws.Range("O2").Value = 19.17549271906796

$ws.Range("B3").Value = 11.22455376164815
$ws.Range("C3").Value = 7.952968402873063
$ws.Range("D3").Value = 3.908799376977902
$ws.Range("E3").Value = 11.68535303170106
$ws.Range("F3").Value = 21.59894531516995
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("M3").Value = 13.95445261415484
$ws.Range("N3").Value = 16.61195351724882
$ws.Range("O3").Value = 19.17683231908308

$ws.Range("B4").Value = 10.86815566977353
$ws.Range("C4").Value = 7.716795311432782
$ws.Range("D4").Value = 3.886432356697304
$ws.Range("E4").Value = 11.61800286406398
$ws.Range("F4").Value = 21.55145764633516
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("M4").Value = 13.78690265558002
$ws.Range("N4").Value = 16.64987138642316
$ws.Range("O4").Value = 19.18361076711308

$ws.Range("B5").Value = 10.71978362971631
$ws.Range("C5").Value = 7.61801927051962
$ws.Range("D5").Value = 3.877214260738068
$ws.Range("E5").Value = 11.59157444985725
$ws.Range("F5").Value = 21.5339975515409
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("M5").Value = 13.71888354482406
$ws.Range("N5").Value = 16.66578571420686
$ws.Range("O5").Value = 19.18786801577225

$ws.Range("B6").Value = 10.69496330934203
$ws.Range("C6").Value = 7.601467156893802
$ws.Range("D6").Value = 3.875677487036029
$ws.Range("E6").Value = 11.58724828360266
$ws.Range("F6").Value = 21.53121293854136
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("M6").Value = 13.70760751819031
$ws.Range("N6").Value = 16.66845624462323
$ws.Range("O6").Value = 19.18866512065622

$ws.Range("B7").Value = 10.86616709788747
$ws.Range("C7").Value = 7.715473326262935
$ws.Range("D7").Value = 3.886308451010633
$ws.Range("E7").Value = 11.61764228522601
$ws.Range("F7").Value = 21.55121449715293
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("M7").Value = 13.7859841487944
$ws.Range("N7").Value = 16.65008413844099
$ws.Range("O7").Value = 19.18366213319532

$ws.Range("B8").Value = 11.59358998256362
$ws.Range("C8").Value = 8.196096751481319
$ws.Range("D8").Value = 3.932301681410655
$ws.Range("E8").Value = 11.7606174827802
$ws.Range("F8").Value = 21.65604328437831
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("M8").Value = 14.13385486837878
$ws.Range("N8").Value = 16.57306912862203
$ws.Range("O8").Value = 19.17471718070536

$ws.Range("B9").Value = 12.9079474125148
$ws.Range("C9").Value = 9.053186629045365
$ws.Range("D9").Value = 4.018912906888967
$ws.Range("E9").Value = 12.0707577005332
$ws.Range("F9").Value = 21.92001595147791
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("M9").Value = 14.81694917484897
$ws.Range("N9").Value = 16.43654429983507
$ws.Range("O9").Value = 19.20450340591524

$ws.Range("B10").Value = 13.79749841492932
$ws.Range("C10").Value = 9.627619568836547
$ws.Range("D10").Value = 4.080121807558141
$ws.Range("E10").Value = 12.31473336545955
$ws.Range("F10").Value = 22.14827884224996
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("M10").Value = 15.3139558856328
$ws.Range("N10").Value = 16.34500352820965
$ws.Range("O10").Value = 19.25526019056269

$ws.Range("B11").Value = 14.18442696672206
$ws.Range("C11").Value = 9.876441002296087
$ws.Range("D11").Value = 4.107391319341644
$ws.Range("E11").Value = 12.42875875494145
$ws.Range("F11").Value = 22.25924311960841
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("M11").Value = 15.53791975109936
$ws.Range("N11").Value = 16.30524626829379
$ws.Range("O11").Value = 19.28459914942769

$ws.Range("B12").Value = 14.32831695990069
$ws.Range("C12").Value = 9.968834400715549
$ws.Range("D12").Value = 4.117630846579054
$ws.Range("E12").Value = 12.4723338514223
$ws.Range("F12").Value = 22.30225552077471
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("M12").Value = 15.62233863622219
$ws.Range("N12").Value = 16.29046114303009
$ws.Range("O12").Value = 19.29660382663612

$ws.Range("B13").Value = 14.2974458584164
$ws.Range("C13").Value = 9.949017686367631
$ws.Range("D13").Value = 4.11542950938904
$ws.Range("E13").Value = 12.46293221333026
$ws.Range("F13").Value = 22.2929484180356
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("M13").Value = 15.60417623693129
$ws.Range("N13").Value = 16.29363338762593
$ws.Range("O13").Value = 19.2939786921497

$ws.Range("B14").Value = 14.19631803800135
$ws.Range("C14").Value = 9.884079107880627
$ws.Range("D14").Value = 4.108235495598493
$ws.Range("E14").Value = 12.43233600130051
$ws.Range("F14").Value = 22.262762073781
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("M14").Value = 15.54487326041683
$ws.Range("N14").Value = 16.30402447978966
$ws.Range("O14").Value = 19.28556887881001

$ws.Range("B15").Value = 14.13402945837205
$ws.Range("C15").Value = 9.844063152631021
$ws.Range("D15").Value = 4.10381752177371
$ws.Range("E15").Value = 12.41364530427748
$ws.Range("F15").Value = 22.24440034516831
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("M15").Value = 15.50849504872
$ws.Range("N15").Value = 16.31042446803561
$ws.Range("O15").Value = 19.28053399446573

$ws.Range("B16").Value = 13.77184759915862
$ws.Range("C16").Value = 9.611104132074695
$ws.Range("D16").Value = 4.078327754471355
$ws.Range("E16").Value = 12.30733921804141
$ws.Range("F16").Value = 22.14116777241254
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("M16").Value = 15.29926910453792
$ws.Range("N16").Value = 16.34763960658738
$ws.Range("O16").Value = 19.25346826139264

$ws.Range("B17").Value = 13.54505633832472
$ws.Range("C17").Value = 9.464966124799165
$ws.Range("D17").Value = 4.062540385874824
$ws.Range("E17").Value = 12.24287348476386
$ws.Range("F17").Value = 22.07964104083324
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("M17").Value = 15.1703086935997
$ws.Range("N17").Value = 16.37095199654326
$ws.Range("O17").Value = 19.23846249190684

$ws.Range("B18").Value = 13.41294679643635
$ws.Range("C18").Value = 9.379737287584121
$ws.Range("D18").Value = 4.053406058119783
$ws.Range("E18").Value = 12.20608310727694
$ws.Range("F18").Value = 22.04492533293595
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("M18").Value = 15.09593870589634
$ws.Range("N18").Value = 16.38453815840131
$ws.Range("O18").Value = 19.23042011826203

$ws.Range("B19").Value = 13.36793349707951
$ws.Range("C19").Value = 9.350679739336108
$ws.Range("D19").Value = 4.050304212440587
$ws.Range("E19").Value = 12.19367732295296
$ws.Range("F19").Value = 22.03328769283292
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("M19").Value = 15.07072743199313
$ws.Range("N19").Value = 16.38916871425659
$ws.Range("O19").Value = 19.22779829391822

$ws.Range("B20").Value = 13.56937163383246
$ws.Range("C20").Value = 9.480644545602029
$ws.Range("D20").Value = 4.064226582327617
$ws.Range("E20").Value = 12.24970639696949
$ws.Range("F20").Value = 22.08612125535534
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("M20").Value = 15.18405760373371
$ws.Range("N20").Value = 16.3684519903201
$ws.Range("O20").Value = 19.23999899416603

$ws.Range("B21").Value = 14.22609372670871
$ws.Range("C21").Value = 9.903203043570047
$ws.Range("D21").Value = 4.110350941735045
$ws.Range("E21").Value = 12.44131242794254
$ws.Range("F21").Value = 22.27160186204832
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("M21").Value = 15.56230326213471
$ws.Range("N21").Value = 16.30096504183578
$ws.Range("O21").Value = 19.28801480254862

$ws.Range("B22").Value = 14.63993320616985
$ws.Range("C22").Value = 10.16869196476059
$ws.Range("D22").Value = 4.139987805146061
$ws.Range("E22").Value = 12.56882783048828
$ws.Range("F22").Value = 22.39859357624989
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("M22").Value = 15.80719114238382
$ws.Range("N22").Value = 16.25843224671676
$ws.Range("O22").Value = 19.32460786242715

$ws.Range("B23").Value = 14.42048868112876
$ws.Range("C23").Value = 10.02798213433335
$ws.Range("D23").Value = 4.124217902045276
$ws.Range("E23").Value = 12.50057459984157
$ws.Range("F23").Value = 22.33029898784199
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("M23").Value = 15.67672895778945
$ws.Range("N23").Value = 16.28098911523044
$ws.Range("O23").Value = 19.30460221587668

$ws.Range("B24").Value = 13.55838405297593
$ws.Range("C24").Value = 9.473560104583754
$ws.Range("D24").Value = 4.063464433109188
$ws.Range("E24").Value = 12.24661638753287
$ws.Range("F24").Value = 22.08318950193505
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("M24").Value = 15.17784242795886
$ws.Range("N24").Value = 16.36958167125508
$ws.Range("O24").Value = 19.23930252003273

$ws.Range("B25").Value = 12.56526752789004
$ws.Range("C25").Value = 8.830831317804964
$ws.Range("D25").Value = 3.995890858450632
$ws.Range("E25").Value = 11.98387268917391
$ws.Range("F25").Value = 21.84247297828357
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("M25").Value = 14.63262776554352
$ws.Range("N25").Value = 16.47193339399939
$ws.Range("O25").Value = 19.19137171584733
